$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to 2-decimal precision (custom accuracy)
$ws.Range("B5").Value = 15.37
$ws.Range("C5").Value = 11.22
$ws.Range("D5").Value = 1.07
$ws.Range("E5").Value = 33.14
$ws.Range("F5").Value = 27.45
$ws.Range("G5").Value = 12.1
$ws.Range("H5").Value = 46.72
$ws.Range("I5").Value = 18.62
$ws.Range("J5").Value = 8.210000000000001
$ws.Range("K5").Value = 12.34
$ws.Range("L5").Value = 13.29
$ws.Range("M5").Value = 13.99
$ws.Range("N5").Value = 3.86
$ws.Range("O5").Value = 12.03
$ws.Range("P5").Value = 17.08
$ws.Range("Q5").Value = 10.18
$ws.Range("R5").Value = 0.83
$ws.Range("S5").Value = 0.68
$ws.Range("T5").Value = 175.8
$ws.Range("U5").Value = 33.65
$ws.Range("V5").Value = 11.1
$ws.Range("W5").Value = 22.55
$ws.Range("X5").Value = 12.02
$ws.Range("Y5").Value = 1.52
$ws.Range("Z5").Value = 22.61
$ws.Range("AA5").Value = 9.81
$ws.Range("AB5").Value = 8.75
$ws.Range("AC5").Value = 10.27
$ws.Range("AD5").Value = 13.97
$ws.Range("AE5").Value = 0.5600000000000001
$ws.Range("AF5").Value = 42.27
$ws.Range("AG5").Value = 6.24
$ws.Range("AH5").Value = 13.88

# Remove row 6 (reduce dataset)
$ws.Rows.Item(6).Delete()
